$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 11 (keep header row 1 and first data row 2)
$ws.Range("A3:F11").EntireRow.Delete()

# Update row 2 with the new job listing data
$ws.Range("A2").Value = "Backend Developer - Laravel"
$ws.Range("B2").Value = "Avrioc Technologies"
$ws.Range("C2").Value = "Abu Dhabi, UAE"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "25 days ago"
$ws.Range("F2").Value = "https://www.gulftalent.com/uae/jobs/backend-developer-laravel-321947"
